# Several updates made to test cases and suites against TDES 21.
#
# The only content-level change in this revision is on the "CV" sheet:
# a new test-case row is appended to the template list ("docgen-pdf-template-fill"),
# filling the previously-blank A35 cell, and the selection follows the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CV")
$ws.Activate()

# New template/test case entry added at the first previously-empty row.
$ws.Range("A35").Value = "docgen-pdf-template-fill"

# Reflect the new row as the active selection (matches the saved view state).
$ws.Range("A35").Select() | Out-Null

# Page setup touched as part of this save (portrait orientation).
$ws.PageSetup.Orientation = 1
